# Update "Datos actualizados" timestamp (A1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 1 de Mayo de 2020 a las 16:52"

# Canada (row 15) - refreshed stats
$ws.Range("B15").Value = 53657
$ws.Range("C15").Value = 421
$ws.Range("D15").Value = 22043
$ws.Range("E15").Value = 28391
$ws.Range("G15").Value = 39
$ws.Range("H15").Value = 3223

# Paises Bajos (row 17) - refreshed stats
$ws.Range("F17").Value = 735

# Pakistan (row 27) - refreshed stats
$ws.Range("B27").Value = 17611
$ws.Range("C27").Value = 1138
$ws.Range("E27").Value = 12890
$ws.Range("G27").Value = 45
$ws.Range("H27").Value = 406

# Barein (row 62) - refreshed stats
$ws.Range("B62").Value = 3169
$ws.Range("C62").Value = 129
$ws.Range("D62").Value = 1553
$ws.Range("E62").Value = 1608

# Togo overtakes Cabo Verde and Camboya in the ranking (sorted by total cases),
# so it moves from row 148 up to row 146; Cabo Verde and Camboya shift down one row.
$ws.Range("A146").Value = "Togo"
$ws.Range("B146").Value = 123
$ws.Range("C146").Value = 7
$ws.Range("D146").Value = 66
$ws.Range("E146").Value = 48
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 9

$ws.Range("A147").Value = "Cabo Verde"
$ws.Range("B147").Value = 122
$ws.Range("C147").Value = 1
$ws.Range("D147").Value = 4
$ws.Range("E147").Value = 117
$ws.Range("F147").Value = 0
$ws.Range("G147").Value = 0
$ws.Range("H147").Value = 1

$ws.Range("A148").Value = "Camboya"
$ws.Range("B148").Value = 122
$ws.Range("C148").Value = 0
$ws.Range("D148").Value = 119
$ws.Range("E148").Value = 3
$ws.Range("F148").Value = 1
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 0
